# update cornstover results after bug fix, exclude weird BMP results from lactic
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 corresponds to biorefinery "cs" (cornstover)
$ws.Range("B5").Value = 2.129246523749677
$ws.Range("C5").Value = 1.96345339780637
$ws.Range("D5").Value = 1.134130246115497
$ws.Range("E5").Value = 2.171047991949717
$ws.Range("F5").Value = 0.07786469255393635
$ws.Range("G5").Value = 0.4673560654131045
$ws.Range("H5").Value = 0.7272457978082166
$ws.Range("I5").Value = -0.3089006207838001
$ws.Range("J5").Value = -2.407826242110738
$ws.Range("K5").Value = 3.23678887278546
$ws.Range("L5").Value = 1.42475408137767
$ws.Range("M5").Value = 4.310883678348473
$ws.Range("O5").Value = 35.47961999108387
$ws.Range("P5").Value = 0.2820027742644197
$ws.Range("R5").Value = 10756.2962108778
$ws.Range("S5").Value = 0.8211167300876979
